# Updated BaseTest / Added IAutoConst / Updated ValidLogin
# Adds a "Password" / "manager" credential pair next to the existing
# "UserName" / "admin" pair on Script1, and selects the resulting A1:B2
# block (mirroring how Excel leaves the range selected after data entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Password"
$ws.Range("B2").Value = "manager"

$ws.Range("A1:B2").Select()
